# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Router labels in column A for rows 45-50 ---
# Old sequence was R1,R1,R2,R2,R3,R3 -> new sequence is R1,R2,R2,R3,R3,R1
$ws.Range("A46").Value = "R2"
$ws.Range("A48").Value = "R3"
$ws.Range("A50").Value = "R1"

# --- Expand the print area to cover the extended table (rows 1-50) ---
$ps = $ws.PageSetup
$ps.PrintArea = '$A$1:$L$50'

# --- Enable "print gridlines" (adds <printOptions gridLines="1"/>) ---
$ps.PrintGridlines = $true

# --- Page margins switched to the metric (cm) equivalents ---
# 1.8 cm / 1.9 cm / 0.8 cm expressed in points (1 cm = 28.3464566929 pt)
$ps.LeftMargin = 51.0236220472441
$ps.RightMargin = 51.0236220472441
$ps.TopMargin = 53.85826771653544
$ps.BottomMargin = 53.85826771653544
$ps.HeaderMargin = 22.677165354330707
$ps.FooterMargin = 22.677165354330707

# --- Move the active selection to P43 (was P44) ---
$ws.Range("P43").Select() | Out-Null
